$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-10 Monday" "2025-03-11 Tuesday"

Replace-Text "814÷4=" "646÷7="
Replace-Text "173÷9=" "442÷6="
Replace-Text "412÷8=" "115÷2="
Replace-Text "624÷9=" "136÷7="
Replace-Text "705÷9=" "798÷4="

Replace-Text "570÷5=" "775÷6="
Replace-Text "520÷2=" "190÷7="
Replace-Text "738÷2=" "620÷2="
Replace-Text "740÷9=" "351÷8="
Replace-Text "879÷3=" "730÷6="

Replace-Text "391÷7=" "611÷4="
Replace-Text "992÷3=" "789÷4="
Replace-Text "290÷4=" "215÷8="
Replace-Text "842÷8=" "956÷4="
Replace-Text "965÷2=" "688÷9="

Replace-Text "586÷5=" "638÷3="
Replace-Text "400÷9=" "757÷6="
Replace-Text "141÷3=" "983÷2="
Replace-Text "764÷3=" "705÷8="
Replace-Text "175÷7=" "636÷2="

Replace-Text "458÷9=" "148÷5="
Replace-Text "284÷8=" "432÷8="
Replace-Text "598÷4=" "376÷2="
Replace-Text "432÷9=" "602÷8="
Replace-Text "287÷2=" "389÷6="
